$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Value = "'.473"
$ws.Range("D32").Value = "'.755"
$ws.Range("E32").Value = 61
$ws.Range("F32").Value = 360
$ws.Range("G32").Value = 108
$ws.Range("H32").Value = 74
$ws.Range("I32").Value = 31
$ws.Range("J32").Value = 17
$ws.Range("K32").Value = 35
$ws.Range("L32").Value = "'16/44"
$ws.Range("C33").Value = "'.472"
$ws.Range("D33").Value = "'.655"
$ws.Range("E33").Value = 18
$ws.Range("F33").Value = 171
$ws.Range("G33").Value = 78
$ws.Range("H33").Value = 55
$ws.Range("I33").Value = 14
$ws.Range("J33").Value = 8
$ws.Range("K33").Value = 25
$ws.Range("L33").Value = "'12/33"
$ws.Range("C34").Value = "'.452"
$ws.Range("D34").Value = "'.814"
$ws.Range("E34").Value = 29
$ws.Range("F34").Value = 196
$ws.Range("G34").Value = 68
$ws.Range("H34").Value = 51
$ws.Range("I34").Value = 16
$ws.Range("J34").Value = 6
$ws.Range("K34").Value = 36
$ws.Range("L34").Value = "'14/38"
$ws.Range("C35").Value = "'.482"
$ws.Range("D35").Value = "'.797"
$ws.Range("E35").Value = 31
$ws.Range("F35").Value = 290
$ws.Range("G35").Value = 82
$ws.Range("H35").Value = 71
$ws.Range("I35").Value = 22
$ws.Range("J35").Value = 13
$ws.Range("K35").Value = 34
$ws.Range("L35").Value = "'15/36"
$ws.Range("C36").Value = "'.535"
$ws.Range("D36").Value = "'.781"
$ws.Range("E36").Value = 26
$ws.Range("F36").Value = 203
$ws.Range("G36").Value = 67
$ws.Range("H36").Value = 36
$ws.Range("I36").Value = 11
$ws.Range("K36").Value = 22
$ws.Range("L36").Value = "'13/39"
$ws.Range("C37").Value = "'.504"
$ws.Range("D37").Value = "'.896"
$ws.Range("E37").Value = 32
$ws.Range("F37").Value = 326
$ws.Range("G37").Value = 110
$ws.Range("H37").Value = 86
$ws.Range("I37").Value = 19
$ws.Range("J37").Value = 16
$ws.Range("K37").Value = 57
$ws.Range("L37").Value = "'17/40"
$ws.Range("C38").Value = "'.504"
$ws.Range("D38").Value = "'.700"
$ws.Range("E38").Value = 25
$ws.Range("F38").Value = 304
$ws.Range("G38").Value = 112
$ws.Range("H38").Value = 40
$ws.Range("I38").Value = 25
$ws.Range("J38").Value = 8
$ws.Range("K38").Value = 33
$ws.Range("L38").Value = "'14/36"
$ws.Range("C39").Value = "'.485"
$ws.Range("D39").Value = "'.836"
$ws.Range("E39").Value = 36
$ws.Range("F39").Value = 339
$ws.Range("G39").Value = 113
$ws.Range("H39").Value = 71
$ws.Range("I39").Value = 25
$ws.Range("J39").Value = 17
$ws.Range("K39").Value = 39
$ws.Range("L39").Value = "'16/36"
$ws.Range("C40").Value = "'.469"
$ws.Range("D40").Value = "'.873"
$ws.Range("E40").Value = 22
$ws.Range("F40").Value = 266
$ws.Range("G40").Value = 116
$ws.Range("H40").Value = 66
$ws.Range("I40").Value = 18
$ws.Range("J40").Value = 14
$ws.Range("K40").Value = 41
$ws.Range("L40").Value = "'15/35"
$ws.Range("D41").Value = "'.793"
$ws.Range("E41").Value = 48
$ws.Range("F41").Value = 375
$ws.Range("G41").Value = 100
$ws.Range("H41").Value = 82
$ws.Range("I41").Value = 29
$ws.Range("J41").Value = 14
$ws.Range("K41").Value = 42
$ws.Range("L41").Value = "'18/43"
$ws.Range("L43").Value = "'0/24"
$ws.Range("L45").Value = "'0/29"
$ws.Range("L48").Value = "'0/26"
$ws.Range("L50").Value = "'0/24"
$ws.Range("L54").Value = "'0/30"
$ws.Range("L57").Value = "'0/28"
$ws.Range("L59").Value = "'0/32"
$ws.Range("L61").Value = "'0/28"
$ws.Range("L64").Value = "'0/29"
$ws.Range("L67").Value = "'0/31"
$ws.Range("L69").Value = "'0/31"
$ws.Range("L70").Value = "'0/32"
$ws.Range("L74").Value = "'0/16"
$ws.Range("L77").Value = "'0/16"
$ws.Range("L78").Value = "'0/18"
$ws.Range("L79").Value = "'0/18"
$ws.Range("L84").Value = "'0/30"
$ws.Range("L86").Value = "'0/31"
$ws.Range("L87").Value = "'0/29"
$ws.Range("L88").Value = "'0/24"
$ws.Range("L94").Value = "'0/27"
$ws.Range("L95").Value = "'0/30"
$ws.Range("L98").Value = "'0/30"
$ws.Range("L99").Value = "'0/27"
$ws.Range("L103").Value = "'0/31"
$ws.Range("L104").Value = "'0/29"
$ws.Range("L105").Value = "'0/29"
$ws.Range("L111").Value = "'0/31"
$ws.Range("L113").Value = "'0/27"
$ws.Range("L114").Value = "'0/28"
$ws.Range("L115").Value = "'0/29"
$ws.Range("L120").Value = "'0/32"
$ws.Range("L123").Value = "'0/25"
$ws.Range("L124").Value = "'0/26"
$ws.Range("L128").Value = "'0/30"
$ws.Range("L130").Value = "'0/25"
$ws.Range("L133").Value = "'0/26"
$ws.Range("L135").Value = "'0/34"
$ws.Range("L138").Value = "'0/26"
$ws.Range("L140").Value = "'0/28"
$ws.Range("L144").Value = "'0/28"
$ws.Range("L147").Value = "'0/30"
$ws.Range("L149").Value = "'0/31"
$ws.Range("L151").Value = "'0/29"
$ws.Range("L154").Value = "'0/35"
$ws.Range("L157").Value = "'0/35"
$ws.Range("L159").Value = "'0/35"
$ws.Range("L160").Value = "'0/38"
$ws.Range("L164").Value = "'0/29"
$ws.Range("L167").Value = "'0/29"
$ws.Range("L168").Value = "'0/33"
$ws.Range("L169").Value = "'0/26"
$ws.Range("L174").Value = "'0/27"
$ws.Range("L176").Value = "'0/32"
$ws.Range("L177").Value = "'0/27"
$ws.Range("L178").Value = "'0/28"
$ws.Range("L184").Value = "'0/31"
$ws.Range("L185").Value = "'0/34"
$ws.Range("L188").Value = "'0/29"
$ws.Range("L189").Value = "'0/31"
$ws.Range("L193").Value = "'0/30"
$ws.Range("L194").Value = "'0/28"
$ws.Range("L195").Value = "'0/30"
$ws.Range("L201").Value = "'0/31"
